{"js": "const body = context.document.body;\n\n// The mapping table's \"Mapeamento (FHIRPath)\" column referenced the root\n// resource type (\"Patient.\") in front of the address element paths. Since\n// the cell already establishes the FHIRPath is relative to Patient, the\n// redundant \"Patient.\" prefix is dropped from each of these FHIRPath\n// expressions.\nconst replacements = [\n  [\"Patient.address.where\", \"address.where\"],\n  [\"Patient.address.city\", \"address.city\"],\n  [\"Patient.address.state\", \"address.state\"],\n];\n\nfor (const [searchText, replaceText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The mapping table's \"Mapeamento (FHIRPath)\" column referenced the root\n# resource type (\"Patient.\") in front of the address element paths. Since\n# the cell already establishes the FHIRPath is relative to Patient, the\n# redundant \"Patient.\" prefix is dropped from each of these FHIRPath\n# expressions.\n$pairs = @(\n    @(\"Patient.address.where\", \"address.where\"),\n    @(\"Patient.address.city\", \"address.city\"),\n    @(\"Patient.address.state\", \"address.state\")\n)\n\nforeach ($pair in $pairs) {\n    $searchText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute([ref]$searchText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$replaceText, [ref]2)\n}\n"}
